$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain string/text updates (non-numeric-looking values)
$ws.Range("D2").Value = "68.315.14"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "3.744.42"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").Value = "3.741.95"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").Value = "  -4.63%  "
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "4.374.04"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "3.754.18"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "68.328.67"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D30").Value = "3.893.71"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -5.08%  "
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("E36").Value = "  --%  "
$ws.Range("D37").Value = "3.700.55"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("E39").Value = "  -6.02%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("E46").Value = "  +10.74%  "
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("E48").Value = "  +3.73%  "
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("E51").Value = "  +0.16%  "

# Numeric-looking price values must be forced to stay text (match inlineStr in target);
# briefly apply a Text number format so Excel does not coerce the string to a number,
# then restore the default style so no stray style attribute is left on the cell.
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "592.65"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "165.98"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.518"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.43"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000259"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "36.02"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "17.85"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.98"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.69"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "464.01"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.87"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000146"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.88"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.05"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.76"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.29"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "29.83"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.15"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.17"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.42"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.301"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "43.30"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "46.54"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.47"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "389.40"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "144.03"
$c.Style = "Normal"
